$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.987640333333333
$ws.Range("H2").Value = 5.962921
$ws.Range("I2").Value = 0.0278174819837782
$ws.Range("J2").Value = 0.0278174819837782
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 19.60726733333334
$ws.Range("N2").Value = 58.82180200000001
$ws.Range("O2").Value = 0.1509859438163708
$ws.Range("P2").Value = 0.1509859438163708
$ws.Range("Q2").Value = 38.97219537818245
$ws.Range("R2").Value = 350.749758403642
$ws.Range("S2").Value = 0.004200048771915643
$ws.Range("T2").Value = 0.004200048771915643
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.987640333333333
$ws.Range("H3").Value = 5.962921
$ws.Range("I3").Value = 0.0278174819837782
$ws.Range("J3").Value = 0.0278174819837782
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 70.710031
$ws.Range("N3").Value = 212.130093
$ws.Range("O3").Value = 0.544503249041223
$ws.Range("P3").Value = 0.544503249041223
$ws.Range("Q3").Value = 140.5461095868503
$ws.Range("R3").Value = 1264.914986281653
$ws.Range("S3").Value = 0.01514670932031292
$ws.Range("T3").Value = 0.01514670932031292
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.987640333333333
$ws.Range("H4").Value = 5.962921
$ws.Range("I4").Value = 0.0278174819837782
$ws.Range("J4").Value = 0.0278174819837782
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 29.95517733333334
$ws.Range("N4").Value = 89.865532
$ws.Range("O4").Value = 0.2306701206736283
$ws.Range("P4").Value = 0.2306701206736284
$ws.Range("Q4").Value = 59.54011865988578
$ws.Range("R4").Value = 535.8610679389719
$ws.Range("S4").Value = 0.0064166619260346
$ws.Range("T4").Value = 0.0064166619260346
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.987640333333333
$ws.Range("H5").Value = 5.962921
$ws.Range("I5").Value = 0.0278174819837782
$ws.Range("J5").Value = 0.0278174819837782
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.589065333333332
$ws.Range("N5").Value = 28.767196
$ws.Range("O5").Value = 0.07384068646877778
$ws.Range("P5").Value = 0.0738406864687778
$ws.Range("Q5").Value = 19.05961301550177
$ws.Range("R5").Value = 171.536517139516
$ws.Range("S5").Value = 0.002054061965515041
$ws.Range("T5").Value = 0.002054061965515041
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 56.81334033333334
$ws.Range("H6").Value = 170.440021
$ws.Range("I6").Value = 0.7951157181995667
$ws.Range("J6").Value = 0.7951157181995667
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.60726733333334
$ws.Range("N6").Value = 58.82180200000001
$ws.Range("O6").Value = 0.1509859438163708
$ws.Range("P6").Value = 0.1509859438163708
$ws.Range("Q6").Value = 1113.954352015316
$ws.Range("R6").Value = 10025.58916813784
$ws.Range("S6").Value = 0.1200512971555931
$ws.Range("T6").Value = 0.1200512971555931
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 56.81334033333334
$ws.Range("H7").Value = 170.440021
$ws.Range("I7").Value = 0.7951157181995667
$ws.Range("J7").Value = 0.7951157181995667
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 70.710031
$ws.Range("N7").Value = 212.130093
$ws.Range("O7").Value = 0.544503249041223
$ws.Range("P7").Value = 0.544503249041223
$ws.Range("Q7").Value = 4017.273056183551
$ws.Range("R7").Value = 36155.45750565195
$ws.Range("S7").Value = 0.4329430919234096
$ws.Range("T7").Value = 0.4329430919234096
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 56.81334033333334
$ws.Range("H8").Value = 170.440021
$ws.Range("I8").Value = 0.7951157181995667
$ws.Range("J8").Value = 0.7951157181995667
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.95517733333334
$ws.Range("N8").Value = 89.865532
$ws.Range("O8").Value = 0.2306701206736283
$ws.Range("P8").Value = 0.2306701206736284
$ws.Range("Q8").Value = 1701.853684584019
$ws.Range("R8").Value = 15316.68316125617
$ws.Range("S8").Value = 0.1834094386665927
$ws.Range("T8").Value = 0.1834094386665927
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 56.81334033333334
$ws.Range("H9").Value = 170.440021
$ws.Range("I9").Value = 0.7951157181995667
$ws.Range("J9").Value = 0.7951157181995667
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.589065333333332
$ws.Range("N9").Value = 28.767196
$ws.Range("O9").Value = 0.07384068646877778
$ws.Range("P9").Value = 0.0738406864687778
$ws.Range("Q9").Value = 544.7868322612351
$ws.Range("R9").Value = 4903.081490351115
$ws.Range("S9").Value = 0.05871189045397127
$ws.Range("T9").Value = 0.05871189045397129
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7501196666666666
$ws.Range("H10").Value = 2.250359
$ws.Range("I10").Value = 0.01049809664416703
$ws.Range("J10").Value = 0.01049809664416703
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.60726733333334
$ws.Range("N10").Value = 58.82180200000001
$ws.Range("O10").Value = 0.1509859438163708
$ws.Range("P10").Value = 0.1509859438163708
$ws.Range("Q10").Value = 14.70779683632422
$ws.Range("R10").Value = 132.370171526918
$ws.Range("S10").Value = 0.001585065030095035
$ws.Range("T10").Value = 0.001585065030095035
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.7501196666666666
$ws.Range("H11").Value = 2.250359
$ws.Range("I11").Value = 0.01049809664416703
$ws.Range("J11").Value = 0.01049809664416703
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 70.710031
$ws.Range("N11").Value = 212.130093
$ws.Range("O11").Value = 0.544503249041223
$ws.Range("P11").Value = 0.544503249041223
$ws.Range("Q11").Value = 53.04098488370966
$ws.Range("R11").Value = 477.3688639533869
$ws.Range("S11").Value = 0.005716247731497709
$ws.Range("T11").Value = 0.00571624773149771
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.7501196666666666
$ws.Range("H12").Value = 2.250359
$ws.Range("I12").Value = 0.01049809664416703
$ws.Range("J12").Value = 0.01049809664416703
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 29.95517733333334
$ws.Range("N12").Value = 89.865532
$ws.Range("O12").Value = 0.2306701206736283
$ws.Range("P12").Value = 0.2306701206736284
$ws.Range("Q12").Value = 22.46996763622089
$ws.Range("R12").Value = 202.229708725988
$ws.Range("S12").Value = 0.002421597219753422
$ws.Range("T12").Value = 0.002421597219753423
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.7501196666666666
$ws.Range("H13").Value = 2.250359
$ws.Range("I13").Value = 0.01049809664416703
$ws.Range("J13").Value = 0.01049809664416703
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.589065333333332
$ws.Range("N13").Value = 28.767196
$ws.Range("O13").Value = 0.07384068646877778
$ws.Range("P13").Value = 0.0738406864687778
$ws.Range("Q13").Value = 7.192946491484888
$ws.Range("R13").Value = 64.73651842336399
$ws.Range("S13").Value = 0.0007751866628208661
$ws.Range("T13").Value = 0.0007751866628208663
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 11.90182033333333
$ws.Range("H14").Value = 35.705461
$ws.Range("I14").Value = 0.166568703172488
$ws.Range("J14").Value = 0.166568703172488
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 19.60726733333334
$ws.Range("N14").Value = 58.82180200000001
$ws.Range("O14").Value = 0.1509859438163708
$ws.Range("P14").Value = 0.1509859438163708
$ws.Range("Q14").Value = 233.3621730289692
$ws.Range("R14").Value = 2100.259557260722
$ws.Range("S14").Value = 0.02514953285876702
$ws.Range("T14").Value = 0.02514953285876702
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 11.90182033333333
$ws.Range("H15").Value = 35.705461
$ws.Range("I15").Value = 0.166568703172488
$ws.Range("J15").Value = 0.166568703172488
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 70.710031
$ws.Range("N15").Value = 212.130093
$ws.Range("O15").Value = 0.544503249041223
$ws.Range("P15").Value = 0.544503249041223
$ws.Range("Q15").Value = 841.5780847264303
$ws.Range("R15").Value = 7574.202762537872
$ws.Range("S15").Value = 0.09069720006600278
$ws.Range("T15").Value = 0.09069720006600278
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 11.90182033333333
$ws.Range("H16").Value = 35.705461
$ws.Range("I16").Value = 0.166568703172488
$ws.Range("J16").Value = 0.166568703172488
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 29.95517733333334
$ws.Range("N16").Value = 89.865532
$ws.Range("O16").Value = 0.2306701206736283
$ws.Range("P16").Value = 0.2306701206736284
$ws.Range("Q16").Value = 356.5211386744725
$ws.Range("R16").Value = 3208.690248070252
$ws.Range("S16").Value = 0.03842242286124758
$ws.Range("T16").Value = 0.03842242286124758
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 11.90182033333333
$ws.Range("H17").Value = 35.705461
$ws.Range("I17").Value = 0.166568703172488
$ws.Range("J17").Value = 0.166568703172488
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.589065333333332
$ws.Range("N17").Value = 28.767196
$ws.Range("O17").Value = 0.07384068646877778
$ws.Range("P17").Value = 0.0738406864687778
$ws.Range("Q17").Value = 114.1273327619284
$ws.Range("R17").Value = 1027.145994857356
$ws.Range("S17").Value = 0.0122995473864706
$ws.Range("T17").Value = 0.0122995473864706
